# Auto-generated Excel COM-interop script implementing the Seraph_Profits scheduled-runner update.
# For each affected row, cells either get a new numeric value or (when the diff removes the
# cell entirely) get cleared so the <c> element itself disappears from the saved XML.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 198.83333
$ws.Range("J4").Value = 208.33333
$ws.Range("L4").Value = 208.33333
$ws.Range("N4").Value = -436.33333

$ws.Range("H6").Value = 202.71428
$ws.Range("I6").Value = 202.71428
$ws.Range("K6").Value = 608.14284
$ws.Range("M6").Value = -496.14284

$ws.Range("H40").Value = 1583.1471
$ws.Range("I40").Value = 1780.6364
$ws.Range("K40").Value = 1780.6364
$ws.Range("M40").Value = -1605.6364

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H92").Value = 1299.2
$ws.Range("I92").Value = 1299.2
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1299.2
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -51.20000000000005
$ws.Range("N92").ClearContents()

$ws.Range("H98").Value = 2101.6
$ws.Range("I98").Value = 2101.6
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2101.6
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -603.5999999999999
$ws.Range("N98").ClearContents()

$ws.Range("H106").Value = 20594.238
$ws.Range("I106").Value = 21323.65
$ws.Range("J106").Value = 6006
$ws.Range("K106").Value = 21323.65
$ws.Range("L106").Value = 6006
$ws.Range("M106").Value = -20692.65
$ws.Range("N106").Value = -7268

$ws.Range("H122").Value = 2101.6
$ws.Range("I122").Value = 2101.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6304.799999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3854.799999999999
$ws.Range("N122").ClearContents()

$ws.Range("H137").Value = 2008
$ws.Range("I137").Value = 1936
$ws.Range("K137").Value = 5808
$ws.Range("M137").Value = -3258


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1091.7273
$ws.Range("I74").Value = 892.7778
$ws.Range("J74").Value = 1987
$ws.Range("K74").Value = 892.7778
$ws.Range("L74").Value = 1987
$ws.Range("M74").Value = -18.77779999999996
$ws.Range("N74").Value = -3735

$ws.Range("H77").Value = 1091.7273
$ws.Range("I77").Value = 892.7778
$ws.Range("J77").Value = 1987
$ws.Range("K77").Value = 4463.889
$ws.Range("L77").Value = 9935
$ws.Range("M77").Value = -95.88900000000012
$ws.Range("N77").Value = -18671


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5296
$ws.Range("I31").Value = 4162
$ws.Range("J31").Value = 5782
$ws.Range("K31").Value = 4162
$ws.Range("L31").Value = 5782
$ws.Range("M31").Value = -3867
$ws.Range("N31").Value = -6372

$ws.Range("H33").Value = 23500750
$ws.Range("I33").Value = 23500750
$ws.Range("K33").Value = 23500750
$ws.Range("M33").Value = -23500371

$ws.Range("H34").Value = 5296
$ws.Range("I34").Value = 4162
$ws.Range("J34").Value = 5782
$ws.Range("K34").Value = 4162
$ws.Range("L34").Value = 5782
$ws.Range("M34").Value = -3960
$ws.Range("N34").Value = -6186

$ws.Range("H107").Value = 1194.7084
$ws.Range("I107").Value = 967.94446
$ws.Range("K107").Value = 967.94446
$ws.Range("M107").Value = 952.05554

$ws.Range("H141").Value = 80000
$ws.Range("J141").Value = 80000
$ws.Range("L141").Value = 80000
$ws.Range("N141").Value = -90360


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 19560.727
$ws.Range("J34").Value = 35383.832
$ws.Range("L34").Value = 106151.496
$ws.Range("N34").Value = -106319.496

$ws.Range("H46").Value = 1428888
$ws.Range("I46").Value = 100
$ws.Range("J46").Value = 2000403.2
$ws.Range("K46").Value = 300
$ws.Range("L46").Value = 6001209.6
$ws.Range("M46").Value = -209
$ws.Range("N46").Value = -6001391.6

$ws.Range("H55").Value = 32178.428
$ws.Range("J55").Value = 44999.8
$ws.Range("L55").Value = 134999.4
$ws.Range("N55").Value = -135353.4

$ws.Range("H58").Value = 583
$ws.Range("I58").Value = 425
$ws.Range("J58").Value = 899
$ws.Range("K58").Value = 1275
$ws.Range("L58").Value = 2697
$ws.Range("M58").Value = -1147
$ws.Range("N58").Value = -2953

$ws.Range("H108").Value = 1980.6666
$ws.Range("I108").Value = 376.8
$ws.Range("K108").Value = 1130.4
$ws.Range("M108").Value = 1749.6

$ws.Range("H113").Value = 797.5
$ws.Range("J113").Value = 795
$ws.Range("L113").Value = 2385
$ws.Range("N113").Value = -6725


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1450
$ws.Range("I22").Value = 1850
$ws.Range("J22").Value = 1050
$ws.Range("K22").Value = 1850
$ws.Range("L22").Value = 1050
$ws.Range("M22").Value = -1555
$ws.Range("N22").Value = -1640

$ws.Range("H27").Value = 1450
$ws.Range("I27").Value = 1850
$ws.Range("J27").Value = 1050
$ws.Range("K27").Value = 1850
$ws.Range("L27").Value = 1050
$ws.Range("M27").Value = -1743
$ws.Range("N27").Value = -1264

$ws.Range("H46").Value = 2562.375
$ws.Range("J46").Value = 2562.375
$ws.Range("L46").Value = 2562.375
$ws.Range("N46").Value = -2938.375

$ws.Range("H68").Value = 4020
$ws.Range("J68").Value = 3001.5
$ws.Range("L68").Value = 3001.5
$ws.Range("N68").Value = -4499.5

$ws.Range("H71").Value = 4020
$ws.Range("J71").Value = 3001.5
$ws.Range("L71").Value = 15007.5
$ws.Range("N71").Value = -22495.5

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 809
$ws.Range("I55").Value = 48
$ws.Range("J55").Value = 1189.5
$ws.Range("K55").Value = 48
$ws.Range("L55").Value = 1189.5
$ws.Range("M55").Value = 229
$ws.Range("N55").Value = -1743.5

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H107").Value = 781.5
$ws.Range("I107").Value = 781.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2344.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -424.5
$ws.Range("N107").ClearContents()

$ws.Range("H122").Value = 2109
$ws.Range("I122").Value = 2173.7
$ws.Range("K122").Value = 6521.099999999999
$ws.Range("M122").Value = -4071.099999999999

$ws.Range("H132").Value = 6833.3335
$ws.Range("I132").Value = 6833.3335
$ws.Range("K132").Value = 20500.0005
$ws.Range("M132").Value = -17970.0005

